# Masters added for Feature with postman api url & body data
# Adds a new "lstfeature" table block (rows 44-49) mirroring the existing
# API-call tables already present on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Values, entered in the same cell order the original author used ---
# (this order drives de-duplication order in xl/sharedStrings.xml)
$ws.Range("B45").Value = "GET 1"
$ws.Range("C45").Value = 100
$ws.Range("D45").Value = "http://localhost:52013/api/Feature/GetFeatures"
$ws.Range("E45").Value = "NA"

$ws.Range("B46").Value = "GET 2"
$ws.Range("C46").Value = 100
$ws.Range("D46").Value = "http://localhost:52013/api/Feature/GetFeature/Id=1"
$ws.Range("E46").Value = "NA"

$ws.Range("B47").Value = "ADD"
$ws.Range("C47").Value = 100
$ws.Range("E47").Value = '{"sFeature":"TestPage","nEnteredBy":0,"nUpdatedBy":0}'
$ws.Range("D47").Value = "http://localhost:52013/api/Feature/AddFeature"

$ws.Range("B48").Value = "EDIT"
$ws.Range("C48").Value = 100
$ws.Range("E48").Value = '{"id":24,"sFeature":"Payments","nEnteredBy":0,"nUpdatedBy":0}'
$ws.Range("D48").Value = "http://localhost:52013/api/Feature/EditFeature/Id=24"

$ws.Range("B49").Value = "DELETE"
$ws.Range("C49").Value = 100
$ws.Range("D49").Value = "http://localhost:52013/api/Feature/DeleteFeature"

$ws.Range("A44").Value = "lstfeature"

$ws.Range("E49").Value = '{"id":24,"sFeature":"Payments","nEnteredBy":0,"nUpdatedBy":0}'

# --- Formatting ---
# D45:D48 get a new style: left/center aligned, wrapped text, built on top of
# the font already used by the other "URL/body" cells (fontId 1 / Arial 9).
$ws.Range("D10").Copy() | Out-Null
$ws.Range("D45:D48").PasteSpecial(-4122) | Out-Null
$fmt = $ws.Range("D45:D48")
$fmt.HorizontalAlignment = -4131
$fmt.VerticalAlignment = -4108
$fmt.WrapText = $true
$fmt.ReadingOrder = 1

# D49 reuses the plain existing "URL/body" style (same as D10, D11, ...).
$ws.Range("D10").Copy() | Out-Null
$ws.Range("D49").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Re-set the values after paste-special so text isn't clobbered by the copy.
$ws.Range("D45").Value = "http://localhost:52013/api/Feature/GetFeatures"
$ws.Range("D46").Value = "http://localhost:52013/api/Feature/GetFeature/Id=1"
$ws.Range("D47").Value = "http://localhost:52013/api/Feature/AddFeature"
$ws.Range("D48").Value = "http://localhost:52013/api/Feature/EditFeature/Id=24"
$ws.Range("D49").Value = "http://localhost:52013/api/Feature/DeleteFeature"

# --- View: move selection to the next empty row beneath the new data ---
$ws.Activate()
$ws.Range("A54").Select() | Out-Null

Write-Output "edit applied"
